$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.621.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4884"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2912"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06686"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.893.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07225"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.008"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6698"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.587.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007916"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.136.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.777"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +29.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.037"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.315"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.875"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.412"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.257"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09024"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.938"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05279"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7366"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.089"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.760"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01821"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.678"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9217"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.075"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4376"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.671"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1345"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.439"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05852"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.759"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3932"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.417"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.46%  "
